# Refresh the rolling production-forecast table.
#
# The sheet holds one row per (date, hourly-interval 1..24) slot in
# columns A:D (A=date serial, B=interval, C=model prediction, D=lookup
# key "DD.MM.YYYY"&interval). Every refresh the window rolls forward:
# the oldest slots are dropped off the top and the same number of new
# slots is appended at the bottom, so each existing row just inherits
# the (date, interval) pair that used to sit 47 rows below it, and the
# last 47 rows extend the sequence with brand-new slots. The Prediction
# column is freshly produced by the forecasting model for the new
# window.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 170

# Freshly produced model output (Prediction, column C) for the new
# rolling window, top to bottom.
$predictions = @(
    0,0,0.447,0.286,0.095,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0.061,0.103,0.173,0.229,0.239,0.218,0.141,0.058,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0.068,0.197,0.313,0.401,0.38,0.286,0.184,0.064,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0.061,0.132,0.199,0.245,0.245,0.218,0.146,0.064,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0.068,0.211,0.341,0.451,0.454,0.35,0.18,0.064,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0.064,0.18,0.307,0.423,0.434,0.358,0.184,0.064,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0.064,0.184,0.393,0.471,0.471,0.434,0.2,0.074,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0.064,0.247,0.443,0.547
)

# New window starts at 02.12.2025, interval 13 and advances by one
# hourly interval per row (interval wraps 24 -> 1 and rolls the date).
$day = Get-Date -Year 2025 -Month 12 -Day 2 -Hour 0 -Minute 0 -Second 0
$interval = 13

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $serial = [int]$day.ToOADate()

    $ws.Cells.Item($r, 1).Value = $serial
    $ws.Cells.Item($r, 2).Value = $interval
    $ws.Cells.Item($r, 3).Value = $predictions[$r - $firstRow]
    $ws.Cells.Item($r, 4).Value = ($day.ToString("dd.MM.yyyy") + $interval.ToString())

    if ($interval -eq 24) {
        $interval = 1
        $day = $day.AddDays(1)
    } else {
        $interval = $interval + 1
    }
}
